$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1410
$ws.Range("J43").Value = 1583.3334
$ws.Range("L43").Value = 1583.3334
$ws.Range("N43").Value = -1721.3334
$ws.Range("H51").Value = 5892.4375
$ws.Range("I51").Value = 4570.2856
$ws.Range("J51").Value = 6920.778
$ws.Range("K51").Value = 4570.2856
$ws.Range("L51").Value = 6920.778
$ws.Range("M51").Value = -4086.2856
$ws.Range("N51").Value = -7888.778
$ws.Range("H80").Value = 649751.5600000001
$ws.Range("I80").Value = 1337244.5
$ws.Range("J80").Value = 452.6111
$ws.Range("K80").Value = 4011733.5
$ws.Range("L80").Value = 1357.8333
$ws.Range("M80").Value = -4010735.5
$ws.Range("N80").Value = -3353.8333
$ws.Range("H83").Value = 649751.5600000001
$ws.Range("I83").Value = 1337244.5
$ws.Range("J83").Value = 452.6111
$ws.Range("K83").Value = 12035200.5
$ws.Range("L83").Value = 4073.4999
$ws.Range("M83").Value = -12030208.5
$ws.Range("N83").Value = -14057.4999
$ws.Range("H132").Value = 16668906
$ws.Range("I132").Value = 18183872
$ws.Range("J132").Value = 4276.4
$ws.Range("K132").Value = 54551616
$ws.Range("L132").Value = 12829.2
$ws.Range("M132").Value = -54549086
$ws.Range("N132").Value = -17889.2
$ws.Range("H138").Value = 2332.422
$ws.Range("I138").Value = 2241.8696
$ws.Range("J138").Value = 2427.0908
$ws.Range("K138").Value = 6725.6088
$ws.Range("L138").Value = 7281.2724
$ws.Range("M138").Value = -1585.6088
$ws.Range("N138").Value = -17561.2724
$ws.Range("H141").Value = 5323.75
$ws.Range("I141").Value = 5673.231
$ws.Range("J141").Value = 3809.3333
$ws.Range("K141").Value = 17019.693
$ws.Range("L141").Value = 11427.9999
$ws.Range("M141").Value = -11839.693
$ws.Range("N141").Value = -21787.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 6298
$ws.Range("I5").Value = 346
$ws.Range("J5").Value = 12250
$ws.Range("K5").Value = 346
$ws.Range("L5").Value = 12250
$ws.Range("M5").Value = -234
$ws.Range("N5").Value = -12474
$ws.Range("H32").Value = 8177.5
$ws.Range("I32").Value = 4917.7554
$ws.Range("J32").Value = 19461.23
$ws.Range("K32").Value = 4917.7554
$ws.Range("L32").Value = 19461.23
$ws.Range("M32").Value = -4630.7554
$ws.Range("N32").Value = -20035.23
$ws.Range("H97").Value = 1076.6154
$ws.Range("I97").Value = 1186.2778
$ws.Range("K97").Value = 1186.2778
$ws.Range("M97").Value = -690.2778000000001
$ws.Range("H122").Value = 2878
$ws.Range("I122").Value = 2314.162
$ws.Range("K122").Value = 6942.485999999999
$ws.Range("M122").Value = -4492.485999999999
$ws.Range("H132").Value = 2871.182
$ws.Range("I132").Value = 1856.8235
$ws.Range("K132").Value = 5570.470499999999
$ws.Range("M132").Value = -3040.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 6298
$ws.Range("I4").Value = 346
$ws.Range("J4").Value = 12250
$ws.Range("K4").Value = 346
$ws.Range("L4").Value = 12250
$ws.Range("M4").Value = -231
$ws.Range("N4").Value = -12480
$ws.Range("H133").Value = 48825
$ws.Range("J133").Value = 48825
$ws.Range("L133").Value = 48825
$ws.Range("N133").Value = -58945
$ws.Range("H134").Value = 2918.628
$ws.Range("I134").Value = 2328.2856
$ws.Range("J134").Value = 5501.375
$ws.Range("K134").Value = 6984.8568
$ws.Range("L134").Value = 16504.125
$ws.Range("M134").Value = -4449.8568
$ws.Range("N134").Value = -21574.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 438135.3
$ws.Range("I58").Value = 1729.2858
$ws.Range("J58").Value = 1116989.1
$ws.Range("K58").Value = 1729.2858
$ws.Range("L58").Value = 1116989.1
$ws.Range("M58").Value = -1526.2858
$ws.Range("N58").Value = -1117395.1
$ws.Range("H62").Value = 44219.332
$ws.Range("I62").Value = 2479.6
$ws.Range("J62").Value = 74033.42999999999
$ws.Range("K62").Value = 2479.6
$ws.Range("L62").Value = 74033.42999999999
$ws.Range("M62").Value = -1855.6
$ws.Range("N62").Value = -75281.42999999999
$ws.Range("H65").Value = 44219.332
$ws.Range("I65").Value = 2479.6
$ws.Range("J65").Value = 74033.42999999999
$ws.Range("K65").Value = 12398
$ws.Range("L65").Value = 370167.15
$ws.Range("M65").Value = -9278
$ws.Range("N65").Value = -376407.15
$ws.Range("H96").Value = 40245.332
$ws.Range("J96").Value = 40245.332
$ws.Range("L96").Value = 40245.332
$ws.Range("N96").Value = -45737.332
$ws.Range("H136").Value = 438135.3
$ws.Range("I136").Value = 1729.2858
$ws.Range("J136").Value = 1116989.1
$ws.Range("K136").Value = 5187.857400000001
$ws.Range("L136").Value = 3350967.3
$ws.Range("M136").Value = -2637.857400000001
$ws.Range("N136").Value = -3356067.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1811.4166
$ws.Range("J23").Value = 2382
$ws.Range("L23").Value = 7146
$ws.Range("N23").Value = -7616
$ws.Range("H62").Value = 9004.666999999999
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 9004.666999999999
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H113").Value = 1279.5
$ws.Range("J113").Value = 1175.2273
$ws.Range("L113").Value = 3525.6819
$ws.Range("N113").Value = -7865.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10045.777
$ws.Range("I126").Value = 15424.889
$ws.Range("K126").Value = 46274.667
$ws.Range("M126").Value = -43804.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 499.75
$ws.Range("J22").Value = 499.66666
$ws.Range("L22").Value = 499.66666
$ws.Range("N22").Value = -1089.66666
$ws.Range("H27").Value = 499.75
$ws.Range("J27").Value = 499.66666
$ws.Range("L27").Value = 499.66666
$ws.Range("N27").Value = -713.66666
$ws.Range("H55").Value = 795.1667
$ws.Range("I55").Value = 434.63635
$ws.Range("J55").Value = 1361.7142
$ws.Range("K55").Value = 434.63635
$ws.Range("L55").Value = 1361.7142
$ws.Range("M55").Value = -261.63635
$ws.Range("N55").Value = -1707.7142
$ws.Range("H93").Value = 19609066
$ws.Range("I93").Value = 22223486
$ws.Range("K93").Value = 22223486
$ws.Range("M93").Value = -22222238
$ws.Range("H122").Value = 3829.1904
$ws.Range("I122").Value = 3732.0625
$ws.Range("K122").Value = 11196.1875
$ws.Range("M122").Value = -8746.1875
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15391992
$ws.Range("I81").Value = 2975
$ws.Range("J81").Value = 22231556
$ws.Range("K81").Value = 5950
$ws.Range("L81").Value = 44463112
$ws.Range("M81").Value = -4889
$ws.Range("N81").Value = -44465234
$ws.Range("H84").Value = 15391992
$ws.Range("I84").Value = 2975
$ws.Range("J84").Value = 22231556
$ws.Range("K84").Value = 29750
$ws.Range("L84").Value = 222315560
$ws.Range("M84").Value = -24446
$ws.Range("N84").Value = -222326168
$ws.Range("H100").Value = 13159538
$ws.Range("I100").Value = 14707631
$ws.Range("J100").Value = 744.5
$ws.Range("K100").Value = 29415262
$ws.Range("L100").Value = 1489
$ws.Range("M100").Value = -29414721
$ws.Range("N100").Value = -2571
$ws.Range("H132").Value = 332667.1
$ws.Range("I132").Value = 458879.4
$ws.Range("K132").Value = 1376638.2
$ws.Range("M132").Value = -1374108.2
$ws.Range("H136").Value = 5725.6
$ws.Range("I136").Value = 5865.4585
$ws.Range("K136").Value = 17596.3755
$ws.Range("M136").Value = -15046.3755
